$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the ".deja.deja.deja" suffix to the header titles that will remain
# after the "Year of Treatment" column (B) is removed, i.e. current columns
# C1:J1 (Occasionally employed ... Total).
for ($col = 3; $col -le 10; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $current = $cell.Value()
    $cell.Value = $current + ".deja.deja.deja"
}

# Remove column B ("Year of Treatment"); this shifts C:J left into B:I,
# shrinking the used range from A1:J35 to A1:I35.
$ws.Columns.Item(2).Delete()
